# Adding Daily Standings to app
#
# The Colors sheet gains a leading "team_id" column (used to join the
# team-color lookup table to the new Daily Standings data). Inserting the
# column shifts the existing Team / Primary Color / Secondary Color /
# Tertiary Color columns one position to the right; their contents and
# per-cell styling move with them automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A:A").Insert()

# Header for the newly inserted column.
$ws.Range("A1").Value2 = "team_id"

# team_id values for each team row (rows 2-13), in the existing row order:
# Anchorage Wheelers, Amarillo Armadillos, Death Valley Scorpions,
# Florida Space Rangers, Kingston Mounties, New York Voyagers,
# Outer Banks Aviators, Providence Crabs, San Antonio Sloths,
# State College Swift Steeds, Utah Railroaders, Vancouver Vandals.
$teamIds = @(16, 19, 10, 24, 17, 1, 2, 3, 9, 18, 25, 11)

for ($i = 0; $i -lt $teamIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $teamIds[$i]
}

# Match the saved selection from the authored workbook.
$ws.Range("C11").Select()
